# Update the handback status report with refreshed generation timestamps
# and corrected priority value, as produced by the handback report generator.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for the
# 9f9e9654-... row (row 3) and the feaa4dea-... row (row 5)
$overview.Range("G3").Value = "2016-09-07 04:23:09"
$overview.Range("G5").Value = "2016-09-07 04:23:09"

# zh-cn sheet: Priority (E) changed from "ht" to "mt"
$zhcn.Range("E3").Value = "mt"
$zhcn.Range("E5").Value = "mt"

# zh-cn sheet: Correspond Handoff Datetime (H)
$zhcn.Range("H3").Value = "2016-09-07 04:22:58"
$zhcn.Range("H5").Value = "2016-09-07 04:22:58"

# zh-cn sheet: Correspond Handback DateTime (K)
$zhcn.Range("K3").Value = "2016-09-07 04:23:30"
$zhcn.Range("K5").Value = "2016-09-07 04:23:30"

# de-de sheet: Priority (E) changed from "ht" to "mt"
$dede.Range("E3").Value = "mt"
$dede.Range("E5").Value = "mt"

# de-de sheet: Correspond Handoff Datetime (H) - shares the same underlying
# timestamp string as the Overview sheet's "Latest HO Xliff Generate Date"
$dede.Range("H3").Value = "2016-09-07 04:23:09"
$dede.Range("H5").Value = "2016-09-07 04:23:09"

# de-de sheet: Correspond Handback DateTime (K)
$dede.Range("K3").Value = "2016-09-07 04:23:38"
$dede.Range("K5").Value = "2016-09-07 04:23:38"
